# Service uptake export: rename the "IP_CODE" column header to "DOB"
# (adding date-of-birth to the export, per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "DOB"
